# Apply updated cryptocurrency price/volume figures to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.126.97'
$ws.Range("E2").Value = '  -0.01%  '
$ws.Range("D3").Value = '1.655.70'
$ws.Range("E3").Value = '  -0.12%  '
$ws.Range("E4").Value = '  -0.21%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '217.93'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.57%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.5264'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +1.89%  '
$ws.Range("E7").Value = '  -0.16%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.2612'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.87%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.06352'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +1.27%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '20.45'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -1.52%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07780'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.72%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '4.513'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +1.94%  '
$ws.Range("D13").Value = '1.667.47'
$ws.Range("E13").Value = '  +0.71%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.5495'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +1.48%  '
$ws.Range("D15").Value = '0.0₅8224'
$ws.Range("E15").Value = '  +1.42%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '65.47'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +1.05%  '
$ws.Range("D17").Value = '26.133.68'
$ws.Range("E17").Value = '  -0.09%  '
$ws.Range("E18").Value = '  -0.24%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '4.594'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -0.51%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '191.13'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -0.30%  '
$ws.Range("E21").Value = '  -0.13%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '6.036'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.35%  '
$ws.Range("E23").Value = '  -0.21%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '141.83'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +1.46%  '
$ws.Range("E25").Value = '  +1.06%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '7.247'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.85%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '16.10'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +0.20%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '1.428'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +1.59%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '0.05901'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -1.16%  '
$ws.Range("E30").Value = '  +0.30%  '
$ws.Range("E31").Value = '  -0.97%  '
$ws.Range("E32").Value = '  +0.30%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '1.593'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -0.48%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.9530'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -1.19%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '2.783'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +0.48%  '
$ws.Range("E36").Value = '  -0.55%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.5709'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +0.30%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.01620'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +1.85%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '5.814'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -2.43%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.8496'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.58%  '
$ws.Range("E41").Value = '  -0.10%  '
$ws.Range("D42").Value = '1.029.64'
$ws.Range("E42").Value = '  +2.31%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '102.72'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +2.28%  '
$ws.Range("D44").Value = '1.799.59'
$ws.Range("E44").Value = '  +0.05%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '57.19'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.84%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -0.40%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.4299'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +2.82%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '1.478'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +2.16%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '7.859'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -1.51%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.05153'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -0.28%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.09702'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.12%  '
